$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: A1:E1 ------------------------------------------------
# A1 used to hold "player1" (shared string 0) together with a running
# total in B1. The sheet is reshaped into a score table: A1 becomes a
# label column header, B1:E1 become the four player-name headers.
$ws.Range("A1").Value = "局本場"
$ws.Range("B1").Value = "player1"
$ws.Range("C1").Value = "player2"
$ws.Range("D1").Value = "player3"
$ws.Range("E1").Value = "player4"

# --- Clear out the old A2:B4 score block (rows reshaped below) -------
$ws.Range("A2:B4").ClearContents()

# --- Per-hand score deltas: rows 2 and 3, columns B:E -----------------
$ws.Range("B2").Value = 8000
$ws.Range("C2").Value = -8000

$ws.Range("B3").Value = -2000
$ws.Range("C3").Value = -4000
$ws.Range("D3").Value = 8000
$ws.Range("E3").Value = -2000

# --- Running totals block: columns H (names) / I (formulas) ----------
$ws.Range("H2").Value = "player1"
$ws.Range("I2").Formula = "=25000+SUM(B:B)"

$ws.Range("H3").Value = "player2"
$ws.Range("I3").Formula = "=25000+SUM(C:C)"

$ws.Range("H4").Value = "player3"
$ws.Range("I4").Formula = "=25000+SUM(D:D)"

$ws.Range("H5").Value = "player4"
$ws.Range("I5").Formula = "=25000+SUM(E:E)"

# --- Restore the selection state left behind in the saved file -------
$ws.Range("E11").Select()
